$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-7 (header row 1 untouched).
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics (as in the original sheet).
$rows = @(
    @(2, "ECs",  "Col9a2", "Mag", "ECs",  2, 0.6666666666666666, 0.159112,            0.477336,           0.09257849183497176, 0.09257849183497174,
      1, 0.3333333333333333, 0.09045733333333333, 0.271372, 0.0574491187297735, 0.0574491187297735,
      0.01439284722133333, 0.129535624992,   0.005318552769250659, 0.005318552769250658),

    @(3, "ECs",  "Col9a2", "Mag", "sCs",  2, 0.6666666666666666, 0.159112,            0.477336,           0.09257849183497176, 0.09257849183497174,
      3, 1,                   1.484107,             4.452321,  0.9425508812702265, 0.9425508812702265,
      0.236139232984,      2.125253096856,  0.08725993906572109, 0.08725993906572108),

    @(4, "FAPs", "Col9a2", "Mag", "ECs",  3, 1,                  1.379163333333333,  4.13749,            0.8024590313369979,  0.8024590313369979,
      1, 0.3333333333333333, 0.09045733333333333, 0.271372, 0.0574491187297735, 0.0574491187297735,
      0.1247554373644444,  1.12279893628,   0.04610056416705823, 0.04610056416705823),

    @(5, "FAPs", "Col9a2", "Mag", "sCs",  3, 1,                  1.379163333333333,  4.13749,            0.8024590313369979,  0.8024590313369979,
      3, 1,                   1.484107,             4.452321,  0.9425508812702265, 0.9425508812702265,
      2.046825957143333,   18.42143361429,  0.7563584671699397, 0.7563584671699397),

    @(6, "sCs",  "Col9a2", "Mag", "ECs",  3, 1,                  0.180396,            0.541188,           0.1049624768280303,  0.1049624768280303,
      1, 0.3333333333333333, 0.09045733333333333, 0.271372, 0.0574491187297735, 0.0574491187297735,
      0.016318141104,      0.146863269936,  0.006030001793464614, 0.006030001793464614),

    @(7, "sCs",  "Col9a2", "Mag", "sCs",  3, 1,                  0.180396,            0.541188,           0.1049624768280303,  0.1049624768280303,
      3, 1,                   1.484107,             4.452321,  0.9425508812702265, 0.9425508812702265,
      0.267726966372,      2.409542697348,  0.09893247503456572, 0.09893247503456572)
)

foreach ($r in $rows) {
    $rowIndex = $r[0]
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $r[$col]
    }
}
